$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Clone formatting into the new cells/rows before filling values ----
# New header cell G1 gets the same (bold) formatting as the rest of row 1.
[void]$ws.Range("A1").Copy()
[void]$ws.Range("G1").PasteSpecial(-4122)

# New rows 7-10 need the same per-column formatting as the existing rows:
#   column D -> Arial 10 style
#   columns E,F -> date ("d-mmm") style
[void]$ws.Range("D2").Copy()
[void]$ws.Range("D7:D10").PasteSpecial(-4122)

[void]$ws.Range("E2:F2").Copy()
[void]$ws.Range("E7:F10").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- Header row (row 1): add new "Estimation" column G ----
$ws.Cells.Item(1, 7).Value = "Estimation"

# ---- Data rows ----
# Row 2 (Carson / Manager) - unchanged text, add Estimation
$ws.Cells.Item(2, 1).Value = "In Progress"
$ws.Cells.Item(2, 2).Value = "Carson"
$ws.Cells.Item(2, 3).Value = "Manager"
$ws.Cells.Item(2, 4).Value = "Turn records of items tracked by employees into a table on the desktop app "
$ws.Cells.Item(2, 5).Value = 43900
$ws.Cells.Item(2, 6).Value = 43901
$ws.Cells.Item(2, 7).Value = "1 day "

# Row 3 (Tristen / Member) - new backlog item text
$ws.Cells.Item(3, 1).Value = "In Progress"
$ws.Cells.Item(3, 2).Value = "Tristen"
$ws.Cells.Item(3, 3).Value = "Member"
$ws.Cells.Item(3, 4).Value = "Implement functionality to remove addresses "
$ws.Cells.Item(3, 5).Value = 43900
$ws.Cells.Item(3, 6).Value = 43901
$ws.Cells.Item(3, 7).Value = "1 day"

# Row 4 (Tristen / Member) - new row
$ws.Cells.Item(4, 1).Value = "In Progress"
$ws.Cells.Item(4, 2).Value = "Tristen"
$ws.Cells.Item(4, 3).Value = "Member"
$ws.Cells.Item(4, 4).Value = "Implement DB call to remove addresses (flag as not in use, but keep in DB for record purpose) "
$ws.Cells.Item(4, 5).Value = 43900
$ws.Cells.Item(4, 6).Value = 43901
$ws.Cells.Item(4, 7).Value = "1 day"

# Row 5 (Luke / Warehouse employee) - unchanged text
$ws.Cells.Item(5, 1).Value = "In Progress"
$ws.Cells.Item(5, 2).Value = "Luke"
$ws.Cells.Item(5, 3).Value = "Warehouse employee "
$ws.Cells.Item(5, 4).Value = "Add ability to flag for damaged returns "
$ws.Cells.Item(5, 5).Value = 43900
$ws.Cells.Item(5, 6).Value = 43901
$ws.Cells.Item(5, 7).Value = "1 day"

# Row 6 (Luke / Warehouse employee) - new row
$ws.Cells.Item(6, 1).Value = "In Progress"
$ws.Cells.Item(6, 2).Value = "Luke"
$ws.Cells.Item(6, 3).Value = "Warehouse employee "
$ws.Cells.Item(6, 4).Value = "Add functionality to DB for flagging rentals "
$ws.Cells.Item(6, 5).Value = 43900
$ws.Cells.Item(6, 6).Value = 43901
$ws.Cells.Item(6, 7).Value = "1 day"

# Row 7 (Luke/Tristen/Carson / Librarian) - unchanged text
$ws.Cells.Item(7, 1).Value = "In Progress"
$ws.Cells.Item(7, 2).Value = "Luke/Tristen/Carson"
$ws.Cells.Item(7, 3).Value = "Librarian"
$ws.Cells.Item(7, 4).Value = "Implement view for librarians  "
$ws.Cells.Item(7, 5).Value = 43900
$ws.Cells.Item(7, 6).Value = 43901
$ws.Cells.Item(7, 7).Value = "2 days"

# Row 8 (Luke/Tristen/Carson / Librarian) - re-cased text
$ws.Cells.Item(8, 1).Value = "In Progress"
$ws.Cells.Item(8, 2).Value = "Luke/Tristen/Carson"
$ws.Cells.Item(8, 3).Value = "Librarian"
$ws.Cells.Item(8, 4).Value = "Implement view for viewing users history/information "
$ws.Cells.Item(8, 5).Value = 43900
$ws.Cells.Item(8, 6).Value = 43901
$ws.Cells.Item(8, 7).Value = "1 day"

# Row 9 (Luke/Tristen/Carson / Librarian) - new row
$ws.Cells.Item(9, 1).Value = "In Progress"
$ws.Cells.Item(9, 2).Value = "Luke/Tristen/Carson"
$ws.Cells.Item(9, 3).Value = "Librarian"
$ws.Cells.Item(9, 4).Value = "Implement ability to filter members by who is overdue "
$ws.Cells.Item(9, 5).Value = 43900
$ws.Cells.Item(9, 6).Value = 43901
$ws.Cells.Item(9, 7).Value = "1 day"

# Row 10 (Luke/Tristen/Carson / Librarian) - new row
$ws.Cells.Item(10, 1).Value = "In Progress"
$ws.Cells.Item(10, 2).Value = "Luke/Tristen/Carson"
$ws.Cells.Item(10, 3).Value = "Librarian"
$ws.Cells.Item(10, 4).Value = "Implement librarian table in the DB for logging in"
$ws.Cells.Item(10, 5).Value = 43900
$ws.Cells.Item(10, 6).Value = 43901
$ws.Cells.Item(10, 7).Value = "1 day"

# ---- Selection matches the post-edit workbook state ----
[void]$ws.Range("H6").Select()
